$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 629.1429000000001
$ws.Range("I29").Value = 629.1429000000001
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 1887.4287
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -1606.4287
$ws.Range("N29").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4107.615
$ws.Range("I40").Value = 1633
$ws.Range("J40").Value = 4850
$ws.Range("K40").Value = 1633
$ws.Range("L40").Value = 4850
$ws.Range("M40").Value = -1458
$ws.Range("N40").Value = -5200

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 2909.0908
$ws.Range("I64").Value = 2950
$ws.Range("J64").Value = 2800
$ws.Range("K64").Value = 2950
$ws.Range("L64").Value = 2800
$ws.Range("M64").Value = -2702
$ws.Range("N64").Value = -3296

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 2909.0908
$ws.Range("I67").Value = 2950
$ws.Range("J67").Value = 2800
$ws.Range("K67").Value = 2950
$ws.Range("L67").Value = 2800
$ws.Range("M67").Value = -2092
$ws.Range("N67").Value = -4516

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3248.5625
$ws.Range("I76").Value = 3198.3333
$ws.Range("J76").Value = 4002
$ws.Range("K76").Value = 3198.3333
$ws.Range("L76").Value = 4002
$ws.Range("M76").Value = -2883.3333
$ws.Range("N76").Value = -4632

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 3248.5625
$ws.Range("I79").Value = 3198.3333
$ws.Range("J79").Value = 4002
$ws.Range("K79").Value = 3198.3333
$ws.Range("L79").Value = 4002
$ws.Range("M79").Value = -2106.3333
$ws.Range("N79").Value = -6186

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 3250
$ws.Range("I106").Value = 3000
$ws.Range("J106").Value = 3500
$ws.Range("K106").Value = 3000
$ws.Range("L106").Value = 3500
$ws.Range("M106").Value = -2369
$ws.Range("N106").Value = -4762

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 723293.7
$ws.Range("I116").Value = 1114011.1
$ws.Range("J116").Value = 20002.4
$ws.Range("K116").Value = 1114011.1
$ws.Range("L116").Value = 20002.4
$ws.Range("M116").Value = -1110569.1
$ws.Range("N116").Value = -26886.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11365.191
$ws.Range("I32").Value = 7782.587
$ws.Range("K32").Value = 7782.587
$ws.Range("M32").Value = -7495.587

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 11545701
$ws.Range("I63").Value = 15392535
$ws.Range("K63").Value = 15392535
$ws.Range("M63").Value = -15391849

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 11545701
$ws.Range("I66").Value = 15392535
$ws.Range("K66").Value = 76962675
$ws.Range("M66").Value = -76959243

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H30").Value = 40000
$ws.Range("J30").Value = 40000
$ws.Range("L30").Value = 40000
$ws.Range("N30").Value = -40250

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2849.4
$ws.Range("I105").Value = 2874.25
$ws.Range("J105").Value = 2750
$ws.Range("K105").Value = 2874.25
$ws.Range("L105").Value = 2750
$ws.Range("M105").Value = -1127.25
$ws.Range("N105").Value = -6244

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 21019
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5234.75
$ws.Range("I31").Value = 3060
$ws.Range("J31").Value = 5959.6665
$ws.Range("K31").Value = 3060
$ws.Range("L31").Value = 5959.6665
$ws.Range("M31").Value = -2765
$ws.Range("N31").Value = -6549.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5234.75
$ws.Range("I34").Value = 3060
$ws.Range("J34").Value = 5959.6665
$ws.Range("K34").Value = 3060
$ws.Range("L34").Value = 5959.6665
$ws.Range("M34").Value = -2858
$ws.Range("N34").Value = -6363.6665

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 143042860
$ws.Range("J37").Value = 143042860
$ws.Range("L37").Value = 429128580
$ws.Range("N37").Value = -429128804

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 77.34999999999999
$ws.Range("I38").Value = 34.583332
$ws.Range("K38").Value = 103.749996
$ws.Range("M38").Value = 243.250004

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("M42").ClearContents()
$ws.Range("N42").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 8929551
$ws.Range("I131").Value = 29412828
$ws.Range("J131").Value = 943.38464
$ws.Range("K131").Value = 88238484
$ws.Range("L131").Value = 2830.15392
$ws.Range("M131").Value = -88233444
$ws.Range("N131").Value = -12910.15392

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 11812.4
$ws.Range("I31").Value = 7265.5
$ws.Range("K31").Value = 7265.5
$ws.Range("M31").Value = -6973.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H37").Value = 11812.4
$ws.Range("I37").Value = 7265.5
$ws.Range("K37").Value = 7265.5
$ws.Range("M37").Value = -6988.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6741.857
$ws.Range("I70").Value = 5962.7856
$ws.Range("J70").Value = 8300
$ws.Range("K70").Value = 5962.7856
$ws.Range("L70").Value = 8300
$ws.Range("M70").Value = -5692.7856
$ws.Range("N70").Value = -8840

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 6741.857
$ws.Range("I73").Value = 5962.7856
$ws.Range("J73").Value = 8300
$ws.Range("K73").Value = 5962.7856
$ws.Range("L73").Value = 8300
$ws.Range("M73").Value = -5026.7856
$ws.Range("N73").Value = -10172

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 20835908
$ws.Range("J80").Value = 2907.1428
$ws.Range("L80").Value = 2907.1428
$ws.Range("N80").Value = -4903.1428

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 20835908
$ws.Range("J83").Value = 2907.1428
$ws.Range("L83").Value = 14535.714
$ws.Range("N83").Value = -24519.714

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2099.7678
$ws.Range("I132").Value = 1152.5454
$ws.Range("J132").Value = 2712.6765
$ws.Range("K132").Value = 3457.6362
$ws.Range("L132").Value = 8138.029500000001
$ws.Range("M132").Value = -927.6361999999999
$ws.Range("N132").Value = -13198.0295

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3794.825
$ws.Range("J122").Value = 7500
$ws.Range("L122").Value = 22500
$ws.Range("N122").Value = -27400

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
